$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the old combined comment in D3 into two parts:
# the tail ("Auto growth disabled...") stays on row 3,
# the new server/drive entry (row 9) gets the head of the sentence.
$ws.Range("D3").Value = " Auto growth disabled. Risk of outage removed"

# Add the new row describing the 10.245.10.37 M: drive assessment.
$ws.Range("A9").Value = "10.245.10.37"
$ws.Range("C9").Value = "M:\"
$ws.Range("D9").Value = "Data pointed to drive with space, and there is and old db that needs to be removed. Its backed up on 10.160 this is pending Mr Jose's approval."

# Copy the fill styling used by the other "risk mitigated" comment cells (e.g. D8) onto D9.
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to the newly added cell, matching the saved view state.
$ws.Range("D9").Select()
